$d = $word.ActiveDocument

# --- Step 1: delete the entire "License Information" paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text -eq "License Information`r") {
        $pp.Range.Delete()
        break
    }
}

# --- Step 2: delete the entire "This PDF version is provided under the same license." paragraph ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text -eq "This PDF version is provided under the same license.`r") {
        $pp.Range.Delete()
        break
    }
}

# --- Step 3: delete the standalone italic "Urie" paragraph (the one between "U" heading and the lone space paragraph) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text -eq "Urie`r") {
        $nextPara = $d.Paragraphs($i + 1)
        if ($nextPara.Range.Text -eq " `r") {
            $pp.Range.Delete()
            break
        }
    }
}

# --- Step 4: update the license/attribution paragraph text ---
# Find the paragraph that still contains the bold "Termes clés (Biblica)" run followed by " (French) is based on"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    $t = $pp.Range.Text
    if ($t.StartsWith(" Termes clés (Biblica) (French)") -or $t.StartsWith("Termes clés (Biblica) (French)")) {
        $targetIndex = $i
        break
    }
}

$pTarget = $d.Paragraphs($targetIndex)
$range = $pTarget.Range
$null = $range.Find.Execute("Termes clés (Biblica)", $true, $false, $false, $false, $false, $true, 1, $false, "Biblica Study Notes (Key Terms)", 2)

$pTarget2 = $d.Paragraphs($targetIndex)
$pStart = $pTarget2.Range.Start
$range2 = $d.Range($pStart + 32, $pTarget2.Range.End)
$bigOld = "(French) is based on: Biblica Bible Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0 license."
$bigNew = "© 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
$null = $range2.Find.Execute($bigOld, $true, $false, $false, $false, $false, $true, 1, $false, $bigNew, 2)
